$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (column C) date value for every existing data
#    row (2..259) from 45177 to 45178.
for ($r = 2; $r -le 259; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}

# 2) Row 259 picks up an explicit row height (matches the rest of the sheet).
$ws.Rows.Item(259).RowHeight = 15

# 3) Append the new record as row 260.
$ws.Cells.Item(260, 1).Value = "A 42025-2023"

$ws.Cells.Item(260, 2).Value = 45177
$ws.Cells.Item(260, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(260, 3).Value = 45178
$ws.Cells.Item(260, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(260, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item(260, 5).Value = "SÄTER"

$ws.Cells.Item(260, 7).Value = 7.1
$ws.Cells.Item(260, 8).Value = 0
$ws.Cells.Item(260, 9).Value = 0
$ws.Cells.Item(260, 10).Value = 0
$ws.Cells.Item(260, 11).Value = 0
$ws.Cells.Item(260, 12).Value = 0
$ws.Cells.Item(260, 13).Value = 0
$ws.Cells.Item(260, 14).Value = 0
$ws.Cells.Item(260, 15).Value = 0
$ws.Cells.Item(260, 16).Value = 0
$ws.Cells.Item(260, 17).Value = 0

$ws.Cells.Item(260, 18).WrapText = $true
